# Add a "Save" column (column H) to the s_vals sheet, mirroring the
# header style used by the existing "sum" column (G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled the same as the other header cells (e.g. G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for the new "Save" column.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
